$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G "Baseline Accuracy" header, matching the style of the
# existing header row (F1) by copying its format.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Baseline Accuracy"

# Baseline Accuracy values for data rows 2-54
$values = @(
    0.6409206595108996,
    0.7032200951643245,
    0.7944745675187194,
    0.8035852605953303,
    0.8227656670724061,
    0.8278558518682454,
    0.8351222750912913,
    0.884954446534617,
    0.896536461215005,
    0.9068274869978975,
    0.9271513407841835,
    0.9273726531666113,
    0.9284054442846077,
    0.9309874220795987,
    0.9346759617867286,
    0.947401423776327,
    0.9475489653646122,
    0.948028475526539,
    0.9495407768064623,
    0.951348161262956,
    0.9568072000295084,
    0.965143299767622,
    0.9659916639002619,
    0.9673195381948286,
    0.9732580871233079,
    0.9767621998450813,
    0.9776105639777212,
    0.978938438272288,
    0.9807458227287816,
    0.9813359890819224,
    0.9819999262292058,
    0.9821843532145623,
    0.9824056655969902,
    0.9842499354505552,
    0.9851720703773376,
    0.9879015897606138,
    0.9897089742171075,
    0.9916639002618863,
    0.9923278374091697,
    0.9924753789974549,
    0.9933606285271661,
    0.9933975139242374,
    0.9937663678949504,
    0.9940983364685921,
    0.9950573567924459,
    0.9952417837778024,
    0.9953524399690162,
    0.995942606322157,
    0.9961270333075135,
    0.9983032717347202,
    0.9987090111025045,
    0.9987458964995758,
    0.9996680314263583
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
